# Auto-generated Excel COM-interop edit script
# Updates cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose new value is purely numeric-looking need an explicit
# text number format first, so Excel/IronCalc keeps them as text
# (matching the original inlineStr cell contents) instead of converting
# them to real numbers.
$textFormatCells = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D14", "D15", "D17",
    "D19", "D20", "D22", "D23", "D25", "D27", "D28", "D29", "D31", "D32",
    "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D42", "D44", "D45",
    "D47", "D48", "D49", "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '27.724.93'
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').Value = '1.758.97'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = '326.73'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.4434'
$ws.Range('E7').Value = '  -2.14%  '
$ws.Range('D8').Value = '0.3735'
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('D9').Value = '45.74'
$ws.Range('E9').Value = '  +2.06%  '
$ws.Range('D10').Value = '0.07778'
$ws.Range('E10').Value = '  +3.18%  '
$ws.Range('D11').Value = '1.129'
$ws.Range('E11').Value = '  -1.65%  '
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('E13').Value = '  -3.31%  '
$ws.Range('D14').Value = '6.210'
$ws.Range('E14').Value = '  -1.47%  '
$ws.Range('D15').Value = '7.379'
$ws.Range('E15').Value = '  -2.68%  '
$ws.Range('D16').Value = '1.758.46'
$ws.Range('E16').Value = '  -2.00%  '
$ws.Range('D17').Value = '91.76'
$ws.Range('E17').Value = '  +13.32%  '
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').Value = '0.06226'
$ws.Range('E19').Value = '  -7.83%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = '6.204'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('D23').Value = '0.5334'
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('D24').Value = '27.756.13'
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('D25').Value = '11.67'
$ws.Range('E25').Value = '  -1.30%  '
$ws.Range('E26').Value = '  -3.97%  '
$ws.Range('D27').Value = '20.88'
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').Value = '153.57'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('D29').Value = '2.377'
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('D30').Value = '1.959.24'
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('D31').Value = '129.27'
$ws.Range('E31').Value = '  -2.86%  '
$ws.Range('D32').Value = '1.217'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').Value = '5.793'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').Value = '0.09283'
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('D35').Value = '3.671'
$ws.Range('E35').Value = '  -9.07%  '
$ws.Range('D36').Value = '12.76'
$ws.Range('E36').Value = '  +5.03%  '
$ws.Range('E37').Value = '  -6.35%  '
$ws.Range('D38').Value = '0.02346'
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('D39').Value = '0.6541'
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('D41').Value = '0.06150'
$ws.Range('E41').Value = '  -2.78%  '
$ws.Range('D42').Value = '1.202'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('E43').Value = '  -3.57%  '
$ws.Range('D44').Value = '1.416'
$ws.Range('E44').Value = '  -4.32%  '
$ws.Range('D45').Value = '1.001'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').Value = '0.6049'
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').Value = '3.763'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('D49').Value = '126.14'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = '1.151'
$ws.Range('E51').Value = '  -1.09%  '
